# Analysis question 4 — append a new answer after the existing
# question 3 / Ksmallest+Klargest paragraph, and relocate the
# "_GoBack" bookmark to the end of the freshly typed paragraph
# (mirrors Word's own behaviour of parking _GoBack at the location
# of the most recent edit).

$d = $word.ActiveDocument

# The paragraph that currently ends the document body ("...Ksmallet
# and Klargest will run in O(k + log n) time.")
$lastParaIndex = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($lastParaIndex)

# Position right before that paragraph's end-of-paragraph mark —
# this is where new text/paragraph breaks must be inserted so they
# land inside/after the existing paragraph rather than shifting the
# following section break.
$insertPos = $p.Range.End - 1
$r = $d.Range($insertPos, $insertPos)

# Helper-less inline pattern: insert each chunk of text, grab the
# range that was just inserted, stamp it with the Times New Roman
# font (and italics where needed), then collapse to the end so the
# next chunk lands after it.
function Type-Chunk([string]$text, [bool]$italic) {
    $chunkStart = $r.Start
    $r.InsertAfter($text)
    $chunkRange = $d.Range($chunkStart, $r.End)
    $chunkRange.Font.Name = "Times New Roman"
    if ($italic) {
        $chunkRange.Font.Italic = $true
    }
    $r.Collapse(0)
}

# --- blank paragraph separating question 3's answer from question 4 ---
# NOTE: InsertParagraphAfter() does not advance $r itself, so re-anchor
# a fresh Range from the freshly-created paragraph each time.
$r.InsertParagraphAfter()
$blankParaIndex = $lastParaIndex + 1
$pBlank = $d.Paragraphs.Item($blankParaIndex)
$blankPos = $pBlank.Range.End - 1
$r = $d.Range($blankPos, $blankPos)

# --- new paragraph: "4. For the Hash Table, ..." ---
$r.InsertParagraphAfter()
$q4ParaIndex = $blankParaIndex + 1
$pQ4 = $d.Paragraphs.Item($q4ParaIndex)
$q4Pos = $pQ4.Range.End - 1
$r = $d.Range($q4Pos, $q4Pos)

Type-Chunk "4. " $false
Type-Chunk "For the Hash Table, the worst case performance for RangeCount is O(" $false
Type-Chunk "h" $true
Type-Chunk ") while the worst case performances for Klargest and Ksmallest are O(" $false
Type-Chunk "h" $true
Type-Chunk " + " $false
Type-Chunk "n" $true
Type-Chunk " log" $false
Type-Chunk " k" $true
Type-Chunk ") where " $false
Type-Chunk "h" $true
Type-Chunk " is the size of the hash table, " $false
Type-Chunk "n" $true
Type-Chunk " is the number of elements in the hash table, and " $false
Type-Chunk "k" $true
Type-Chunk " is the number of smallest/largest elements requested in the function call (e.g. when getting the 5 smallest elements, " $false
Type-Chunk "k" $true
Type-Chunk " is 5)." $false

# The "_GoBack" bookmark tracks the most recent edit location; move
# it from its old spot (mid-sentence in the paragraph above) to the
# end of the text we just typed, mirroring what Word does whenever
# you type/insert content.
$goBackPos = $r.Start
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
